$wb = $excel.ActiveWorkbook

# --- Power sheet: update primary 300W input voltage from 220 to 150 ---
$wsPower = $wb.Worksheets.Item("Power")
$wsPower.Range("D4").Value = 150

# --- Sheet2: add Vin/Vout/Turn-ratio/Inductance simulation table ---
# Values are entered in the same order the author typed them so that the
# shared-string table grows in the expected sequence.
$wsSheet2 = $wb.Worksheets.Item("Sheet2")

$wsSheet2.Range("C7").Value = "Vin"
$wsSheet2.Range("E7").Value = "V"
$wsSheet2.Range("C8").Value = "Vout"
$wsSheet2.Range("F7").Value = "T1"
$wsSheet2.Range("F8").Value = "T2"
$wsSheet2.Range("H7").Value = "L1"
$wsSheet2.Range("H8").Value = "L2"
$wsSheet2.Range("J7").Value = "uH"

$wsSheet2.Range("D7").Value = 150
$wsSheet2.Range("G7").Value = 1
$wsSheet2.Range("I7").Value = 600

$wsSheet2.Range("D8").Formula = "=D7*(G8/G7)"
$wsSheet2.Range("E8").Value = "V"
$wsSheet2.Range("G8").Value = 5
$wsSheet2.Range("I8").Formula = "=D8^2*I7/(D7^2)"

# --- Update the selected cell on each sheet ---
$wsPower.Range("D5").Select()
$wsSnubber = $wb.Worksheets.Item("Snubber Review")
$wsSnubber.Range("N4").Select()
$wsSheet2.Range("J19").Select()
$wsSheet2.Activate()
